$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.887.45"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.928.88"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.70"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.537"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.40"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.138"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0837"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").Value = "3.402.29"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.91"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.35"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "2.920.87"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.975"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").Value = "50.854.19"
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.18"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -8.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.09"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.80%  "
$ws.Range("D22").Value = "0.0₃0948"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.87"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.88"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.20"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +9.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.167"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.112"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.56"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.76"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.89"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.32"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.04"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0450"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.97"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.57"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.56"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.29%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.77"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.61"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.20"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.12%  "
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.271"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.89%  "
$ws.Range("D48").Value = "1.999.60"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.20"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0342"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.475"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +11.10%  "
